$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.319.01"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "3.249.62"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'594.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'140.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.242.82"
$ws.Range("E8").Value = "  +3.72%  "
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").Value = "'0.465"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").Value = "'34.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "3.784.48"
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("D17").Value = "3.250.35"
$ws.Range("E17").Value = "  +3.80%  "
$ws.Range("D18").Value = "63.374.78"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "'6.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'473.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").Value = "'14.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("E22").Value = "  +3.93%  "
$ws.Range("D23").Value = "'7.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.55%  "
$ws.Range("D24").Value = "'83.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.45%  "
$ws.Range("D25").Value = "'13.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "'7.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.04%  "
$ws.Range("D29").Value = "'8.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "'2.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.32%  "
$ws.Range("D31").Value = "'27.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").Value = "'2.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").Value = "'5.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").Value = "'52.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "0.0₃0711"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("D39").Value = "'0.0392"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'419.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").Value = "2.986.24"
$ws.Range("E41").Value = "  +3.52%  "
$ws.Range("D42").Value = "'8.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").Value = "'2.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("E44").Value = "  -7.72%  "
$ws.Range("D45").Value = "'0.265"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.92%  "
$ws.Range("D46").Value = "'2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "'25.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "'122.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.97%  "
